# "Initial check-in of translations changes."
#
# The survey/settings sheets use shared-string header labels that follow the
# XLSForm "display.<kind>" convention. This commit renames the generic
# "display.hint" / "display.title" labels that used to sit in the "survey"
# and "settings" sheets to the more specific "display.hint.text" /
# "display.title.text" forms, and likewise promotes the "survey" sheet's
# "display.text" column header to "display.prompt.text" (the "choices"
# sheet's own "display.text" header is untouched - it's a different column).
#
# It also leaves the workbook with the "settings" sheet active/selected
# (previously "survey" was active), with fresh cell selections on the
# sheets that were touched.

$wb = $excel.ActiveWorkbook

$wsSurvey   = $wb.Worksheets.Item("survey")
$wsChoices  = $wb.Worksheets.Item("choices")
$wsModel    = $wb.Worksheets.Item("model")
$wsSettings = $wb.Worksheets.Item("settings")

# --- Rename the translation-column headers on "survey" ---------------------
# D1: display.text  -> display.prompt.text
# E1: display.hint  -> display.hint.text
$wsSurvey.Range("D1").Value = "display.prompt.text"
$wsSurvey.Range("E1").Value = "display.hint.text"

# --- Rename the translation-column header on "settings" --------------------
# C1: display.title -> display.title.text
$wsSettings.Range("C1").Value = "display.title.text"

# --- Update the per-sheet selections ---------------------------------------
# "survey" keeps its own remembered selection, now one cell to the right.
[void]$wsSurvey.Range("E2").Select()

# "choices" and "model" selections are unchanged (C2 / C22 respectively) -
# nothing to do there.

# "settings" becomes the active sheet/tab, selection moves to C2.
[void]$wsSettings.Activate()
[void]$wsSettings.Range("C2").Select()
